$wb = $excel.ActiveWorkbook

# Mapping of sheet name -> list of (row, newValue) updates to column F ("想去人数")
$updates = @{
  "展览" = @(
    @{ Row = 2; Value = 28 }
    @{ Row = 3; Value = 938 }
    @{ Row = 4; Value = 1214 }
    @{ Row = 5; Value = 1654 }
    @{ Row = 6; Value = 884 }
    @{ Row = 7; Value = 545 }
    @{ Row = 8; Value = 2220 }
    @{ Row = 10; Value = 543 }
    @{ Row = 11; Value = 549 }
    @{ Row = 12; Value = 637 }
    @{ Row = 13; Value = 296 }
    @{ Row = 14; Value = 157 }
    @{ Row = 16; Value = 2084 }
    @{ Row = 18; Value = 665 }
    @{ Row = 19; Value = 2580 }
    @{ Row = 21; Value = 18 }
    @{ Row = 24; Value = 312 }
    @{ Row = 25; Value = 1704 }
    @{ Row = 26; Value = 8 }
    @{ Row = 28; Value = 1 }
    @{ Row = 29; Value = 541 }
    @{ Row = 31; Value = 4473 }
    @{ Row = 32; Value = 11 }
  )
  "演出" = @(
    @{ Row = 5; Value = 4184 }
    @{ Row = 7; Value = 39 }
    @{ Row = 8; Value = 54 }
    @{ Row = 11; Value = 57 }
    @{ Row = 12; Value = 8 }
    @{ Row = 14; Value = 308 }
    @{ Row = 23; Value = 1756 }
    @{ Row = 24; Value = 18 }
    @{ Row = 25; Value = 232 }
    @{ Row = 32; Value = 21 }
    @{ Row = 34; Value = 57 }
  )
  "本地生活" = @(
    @{ Row = 4; Value = 1397 }
    @{ Row = 5; Value = 1774 }
    @{ Row = 7; Value = 476 }
    @{ Row = 8; Value = 117 }
  )
  "全部类型" = @(
    @{ Row = 3; Value = 1397 }
    @{ Row = 4; Value = 1774 }
    @{ Row = 5; Value = 476 }
    @{ Row = 9; Value = 28 }
    @{ Row = 10; Value = 938 }
    @{ Row = 11; Value = 1214 }
    @{ Row = 12; Value = 1654 }
    @{ Row = 13; Value = 39 }
    @{ Row = 14; Value = 54 }
    @{ Row = 16; Value = 884 }
    @{ Row = 17; Value = 545 }
    @{ Row = 18; Value = 2220 }
    @{ Row = 20; Value = 543 }
    @{ Row = 21; Value = 549 }
    @{ Row = 22; Value = 637 }
    @{ Row = 23; Value = 296 }
    @{ Row = 24; Value = 57 }
    @{ Row = 25; Value = 157 }
    @{ Row = 27; Value = 308 }
    @{ Row = 29; Value = 2084 }
    @{ Row = 31; Value = 665 }
    @{ Row = 34; Value = 2580 }
    @{ Row = 37; Value = 18 }
    @{ Row = 39; Value = 117 }
    @{ Row = 40; Value = 1756 }
    @{ Row = 41; Value = 1704 }
    @{ Row = 42; Value = 18 }
    @{ Row = 44; Value = 541 }
    @{ Row = 47; Value = 4473 }
    @{ Row = 48; Value = 21 }
    @{ Row = 49; Value = 57 }
  )
}

foreach ($sheetName in $updates.Keys) {
  $ws = $wb.Worksheets.Item($sheetName)
  foreach ($u in $updates[$sheetName]) {
    $ws.Cells.Item($u.Row, 6).Value = $u.Value
  }
}
